$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10
$ws.Range("A10").Value = "B. Wrath"
$ws.Range("B10").Value = "https://codeforces.com/problemset/problem/892/B"
$ws.Range("D10").Value = "last one previous "

# Row 11
$ws.Range("A11").Value = "Sphere(Array-Sub - SubArrays)"
$ws.Range("B11").Value = "https://www.spoj.com/problems/ARRAYSUB/"
$ws.Range("D11").Value = "max element in a fixed window"

# Row 12
$ws.Range("A12").Value = "HOTELS - Hotels Along the Croatian Coast"
$ws.Range("B12").Value = "https://www.spoj.com/problems/HOTELS/"
$ws.Range("D12").Value = "max window sum less than K"

# Row 13 (note insertion order: B and D first, then A, to match shared-string append order)
$ws.Range("B13").Value = "https://codeforces.com/problemset/problem/701/C"
$ws.Range("D13").Value = "Min_Sub_String_Contain_Particular_Sub_String"
$ws.Range("A13").Value = "C. They Are Everywhere"

# Hyperlinks for the two new URL cells that get the Hyperlink style (B10, B11)
$ws.Hyperlinks.Add($ws.Range("B10"), "https://codeforces.com/problemset/problem/892/B")
$ws.Hyperlinks.Add($ws.Range("B11"), "https://www.spoj.com/problems/ARRAYSUB/")

# Column A width
$ws.Columns.Item(1).ColumnWidth = 37.14

# Selection
$ws.Range("A13").Select()
